$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the "New dish" screen glossary entry (B9): field names are now quoted.
$ws.Range("B9").Value = "1. ""Новое блюдо"" screen`n2. Back button `n3. [Сохранить] button`n4. ""Name dish"" field`n5. ""Content info"" field`n6. Segment control`n7. Selected segment`n8. [Выбрать фото] button"

# Rewrite the "New receipt" screen glossary entry (M11): field names are now quoted.
$ws.Range("M11").Value = "1. [Отмена] button `n2. [Добавить] button `n3. ""Name"" receipt field`n4. Category Picker`n5. ""Contetn info"" field`n6. ""Link"" field`n7. Stepper`n8. ""Ingredients"" field`n9. Proportions stepper"

# Update the scroll position / active selection of the sheet view.
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M12").Select()
